$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 8 new rows right before the current last data row (row 24) so
#    that it (and everything below it, i.e. the two footer rows) shift down
#    by 8 rows. This turns:
#       data rows 16-24 (9 rows), footer rows 29-30
#    into:
#       data rows 16-32 (17 rows), footer rows 37-38
# ---------------------------------------------------------------------------
$insertRange = $ws.Range("B24:J31")
$insertRange.Insert(-4121) | Out-Null   # xlShiftDown

# Give the freshly inserted rows the same formatting (borders/fonts/number
# formats) as a normal data row (row 23) so they match the rest of the table
# instead of staying blank-styled.
$ws.Range("B23:J23").Copy() | Out-Null
$ws.Range("B24:J31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Rewrite the whole data table (rows 16-32) with the refreshed data.
# ---------------------------------------------------------------------------
$rowsData = @(
    ,@(16, "CC", "1143359392", "RUBEN DARIO CARAZO SEQUEA", "2111", 1333, 1000000)
    ,@(17, "CC", "1022931164", "MELBA PATRICIA ZARATE GONZALEZ", "2505", 2133, 1600000)
    ,@(18, "CC", "1022931164", "MELBA PATRICIA ZARATE GONZALEZ", "2503", 2133, 1600000)
    ,@(19, "CC", "1023165078", "CARMEN EVITA ANGEL MARTON", "2003", 315000, 11250000)
    ,@(20, "CC", "9294722", "VICTOR MANUEL BENITEZ MONTIEL", "1709", 29509, 737717)
    ,@(21, "CC", "9294722", "VICTOR MANUEL BENITEZ MONTIEL", "1708", 29509, 737717)
    ,@(22, "CC", "9294722", "VICTOR MANUEL BENITEZ MONTIEL", "1707", 29509, 737717)
    ,@(23, "CC", "9294722", "VICTOR MANUEL BENITEZ MONTIEL", "1704", 29509, 737717)
    ,@(24, "CE", "668215", "JOSE RAMON LEZAMA DIAZ", "1704", 29509, 737717)
    ,@(25, "CE", "668215", "JOSE RAMON LEZAMA DIAZ", "1703", 7869, 737717)
    ,@(26, "CC", "1127585376", "MILEIS DE JESUS CONEO ALVAREZ", "1705", 15738, 1520000)
    ,@(27, "CC", "80715230", "JONATHAN CARDOZA LOZANO", "2507", 44000, 1100000)
    ,@(28, "CC", "80715230", "JONATHAN CARDOZA LOZANO", "2506", 44000, 1100000)
    ,@(29, "CC", "80715230", "JONATHAN CARDOZA LOZANO", "2505", 44000, 1100000)
    ,@(30, "CC", "80715230", "JONATHAN CARDOZA LOZANO", "2504", 44000, 1100000)
    ,@(31, "CC", "80715230", "JONATHAN CARDOZA LOZANO", "2503", 44000, 1100000)
    ,@(32, "CC", "80715230", "JONATHAN CARDOZA LOZANO", "2502", 44000, 1100000)
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Range("B$r").Value2 = $row[1]
    $ws.Range("C$r").Value2 = $row[2]
    $ws.Range("D$r").Value2 = $row[3]
    $ws.Range("E$r").Value2 = $row[4]
    $ws.Range("F$r").Value2 = $row[5]
    $ws.Range("G$r").Value2 = $row[6]
}

# ---------------------------------------------------------------------------
# 3) Update the summary fields above the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 755751   # VALOR MORA (total)
$ws.Range("C13").Value2 = 7        # Cant. Trabajadores
$ws.Range("F13").Value2 = 14       # Cant. Periodos

$wb.Save()
